$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-04 19:49:29'
$ws.Range("I2").Value = '1.5 mm'
$ws.Range("E3").Value = '2026-02-04 19:49:31'
$ws.Range("O3").Value = '-4.6 °C'
$ws.Range("E4").Value = '2026-02-04 19:49:34'
$ws.Range("H4").Value = "'80%"
$ws.Range("O4").Value = '6.3 °C'
$ws.Range("E5").Value = '2026-02-04 19:49:37'
$ws.Range("E6").Value = '2026-02-04 19:49:39'
$ws.Range("H6").Value = "'65%"
$ws.Range("E7").Value = '2026-02-04 19:49:41'
$ws.Range("J7").Value = '993.6 hPa'
$ws.Range("E8").Value = '2026-02-04 19:49:44'
$ws.Range("K8").Value = '8.4 MJ/m2'
$ws.Range("O8").Value = '7.5 °C'
$ws.Range("E9").Value = '2026-02-04 19:49:46'
$ws.Range("H9").Value = "'89%"
$ws.Range("E10").Value = '2026-02-04 19:49:49'
$ws.Range("E11").Value = '2026-02-04 19:49:51'
$ws.Range("E12").Value = '2026-02-04 19:49:54'
$ws.Range("H12").Value = "'81%"
$ws.Range("E13").Value = '2026-02-04 19:49:56'
$ws.Range("E14").Value = '2026-02-04 19:49:58'
$ws.Range("O14").Value = '-6.0 °C'
$ws.Range("E15").Value = '2026-02-04 19:50:01'
$ws.Range("H15").Value = "'82%"
$ws.Range("E16").Value = '2026-02-04 19:50:03'
$ws.Range("E17").Value = '2026-02-04 19:50:06'
$ws.Range("O17").Value = '3.3 °C'
$ws.Range("E18").Value = '2026-02-04 19:50:08'
$ws.Range("E19").Value = '2026-02-04 19:50:11'
$ws.Range("H19").Value = "'86%"
$ws.Range("J19").Value = '994.4 hPa'
$ws.Range("O19").Value = '6.9 °C'
$ws.Range("E20").Value = '2026-02-04 19:50:13'
$ws.Range("E21").Value = '2026-02-04 19:50:16'
$ws.Range("J21").Value = '993.0 hPa'
$ws.Range("E22").Value = '2026-02-04 19:50:18'
$ws.Range("O22").Value = '8.3 °C'
$ws.Range("E23").Value = '2026-02-04 19:50:21'
$ws.Range("E24").Value = '2026-02-04 19:50:23'
$ws.Range("H24").Value = "'70%"
$ws.Range("E25").Value = '2026-02-04 19:50:26'
$ws.Range("E26").Value = '2026-02-04 19:50:28'
$ws.Range("E27").Value = '2026-02-04 19:50:30'
$ws.Range("E28").Value = '2026-02-04 19:50:33'
$ws.Range("E29").Value = '2026-02-04 19:50:35'
$ws.Range("E30").Value = '2026-02-04 19:50:38'
$ws.Range("O30").Value = '-5.2 °C'
$ws.Range("E31").Value = '2026-02-04 19:50:40'
$ws.Range("E32").Value = '2026-02-04 19:50:43'
$ws.Range("E33").Value = '2026-02-04 19:50:45'
$ws.Range("E34").Value = '2026-02-04 19:50:48'
$ws.Range("H34").Value = "'87%"
$ws.Range("K34").Value = '7.0 MJ/m2'
$ws.Range("O34").Value = '3.7 °C'
$ws.Range("E35").Value = '2026-02-04 19:50:50'
$ws.Range("I35").Value = '0.7 mm'
$ws.Range("E36").Value = '2026-02-04 19:50:52'
$ws.Range("O36").Value = '7.3 °C'
